$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list (price + 1h volume change) per upstream diff.
# Some Price-column values are numeric-looking strings (e.g. "320.93") that
# must stay plain TEXT (matching the source inlineStr cells), so force the
# cell to Text format before assigning, to stop Excel auto-converting them
# to numbers (which would also corrupt values like "239.00" -> 239).

$ws.Range('D2').Value = '43.165.05'
$ws.Range('E2').Value = '  -5.13%  '
$ws.Range('D3').Value = '2.235.36'
$ws.Range('E3').Value = '  -6.03%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '320.93'
$ws.Range('E5').Value = '  +0.55%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '100.57'
$ws.Range('E6').Value = '  -8.30%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.584'
$ws.Range('E7').Value = '  -8.52%  '
$ws.Range('E8').Value = '  -0.17%  '
$ws.Range('E9').Value = '  -8.46%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '37.18'
$ws.Range('E10').Value = '  -9.85%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '54.44'
$ws.Range('E11').Value = '  -3.02%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0833'
$ws.Range('E12').Value = '  -9.57%  '
$ws.Range('E13').Value = '  -9.99%  '
$ws.Range('E14').Value = '  -1.21%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.871'
$ws.Range('E15').Value = '  -11.82%  '
$ws.Range('D16').Value = '2.573.49'
$ws.Range('E16').Value = '  -6.12%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '14.48'
$ws.Range('E17').Value = '  -6.46%  '
$ws.Range('D18').Value = '2.235.22'
$ws.Range('E18').Value = '  -5.65%  '
$ws.Range('D19').Value = '43.085.85'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.62'
$ws.Range('E20').Value = '  -4.68%  '
$ws.Range('E21').Value = '  -8.88%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.57'
$ws.Range('E22').Value = '  -10.40%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '65.69'
$ws.Range('E23').Value = '  -10.77%  '
$ws.Range('E24').Value = '  -13.32%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '239.00'
$ws.Range('E25').Value = '  -10.45%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.18'
$ws.Range('E26').Value = '  -7.70%  '
$ws.Range('E27').Value = '  -0.06%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '4.04'
$ws.Range('E28').Value = '  +1.15%  '
$ws.Range('E29').Value = '  -1.79%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '10.07'
$ws.Range('E30').Value = '  -10.35%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.38'
$ws.Range('E31').Value = '  -15.58%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '36.12'
$ws.Range('E32').Value = '  -3.30%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '20.45'
$ws.Range('E33').Value = '  -9.47%  '
$ws.Range('E34').Value = '  -7.47%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '153.99'
$ws.Range('E35').Value = '  -8.68%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.68'
$ws.Range('E36').Value = '  -3.25%  '
$ws.Range('E37').Value = '  +7.84%  '
$ws.Range('E38').Value = '  +1.54%  '
$ws.Range('E39').Value = '  -7.52%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '4.47'
$ws.Range('E40').Value = '  -5.53%  '
$ws.Range('E41').Value = '  -10.81%  '
$ws.Range('E42').Value = '  -8.16%  '
$ws.Range('E43').Value = '  -8.30%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.63'
$ws.Range('E44').Value = '  +4.94%  '
$ws.Range('E45').Value = '  -0.06%  '
$ws.Range('D46').Value = '1.738.72'
$ws.Range('E46').Value = '  -7.62%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '86.22'
$ws.Range('E47').Value = '  -12.84%  '
$ws.Range('E48').Value = '  -9.71%  '
$ws.Range('E49').Value = '  -10.41%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '75.67'
$ws.Range('E50').Value = '  -9.97%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '59.07'
$ws.Range('E51').Value = '  -16.29%  '
